$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetups")

# Renumber the meetup blog links so they match the chapter numbers.
$ws.Range("F7").Value  = "[Recording and slides](/blog/meetup_05a_foundation_for_inference/)"
$ws.Range("F8").Value  = "[Recording and slides](/blog/meetup_05b_foundation_for_inference2/)"
$ws.Range("F9").Value  = "[Recording and slides](/blog/meetup_06_inference_for_categorical_data/)"
$ws.Range("F10").Value = "[Recording and slides](/blog/meetup_07_inference_for_numerical_data/)"

# Move the active selection to match the saved view in the source file.
$ws.Activate()
$ws.Range("F11").Select()
